$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring in the same header style/formatting used by the existing header cells
# (bold font, border, centered/top alignment) by copying from H1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
